$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Insert a new row above row 13 - this shifts the existing "Good Drivers"
# rows (13-18) down to (14-19), preserving their per-cell formatting.
$ws.Rows(13).Insert()

# Copy the (now shifted) row 14 formatting down into the new row 13 so the
# freshly inserted driver row matches the table's look (right-aligned
# numbers, thousands separator on the count column, etc.)
$ws.Range("B14:E14").Copy()
$ws.Range("B13:E13").PasteSpecial(-4122)

# New driver entry added at the top of the "Good Drivers" list
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Range("B13").Value = 11128
$ws.Range("D13").Value = 100
# No vintage date on record yet for this freshly-observed driver
$ws.Range("E13").Value = 0

# Weekly refresh: updated client counts for the drivers that shifted down
$ws.Range("B14").Value = 486214
$ws.Range("B15").Value = 79953
$ws.Range("B16").Value = 35355
$ws.Range("B17").Value = 65425
$ws.Range("B18").Value = 117653
